$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> D value (price), only set when changed (blank = unchanged / skip)
$dValues = @{
    2  = "43.185.86"
    3  = "2.340.33"
    4  = "1.00"
    5  = "302.58"
    6  = "95.50"
    10 = "34.16"
    11 = "0.0785"
    12 = "18.70"
    15 = "2.706.58"
    16 = "2.353.99"
    18 = "43.115.95"
    19 = "12.24"
    20 = "6.21"
    21 = "0.0₃0891"
    22 = "68.02"
    23 = "235.93"
    24 = "2.21"
    26 = "2.42"
    27 = "24.63"
    28 = "2.36"
    29 = "9.20"
    30 = "31.56"
    31 = "1.00"
    33 = "0.0736"
    34 = "17.27"
    36 = "4.36"
    39 = "2.75"
    40 = "22.28"
    42 = "117.67"
    43 = "1.938.84"
    45 = "10.04"
    48 = "2.570.37"
    49 = "53.22"
    51 = "72.09"
}

# Map of row -> E value (Volume(1h)); padded with two leading/trailing spaces
$eValues = @{
    2  = "+1.09%"
    3  = "+1.69%"
    4  = "+0.03%"
    5  = "+0.49%"
    6  = "-0.16%"
    7  = "+0.17%"
    8  = "-0.08%"
    9  = "+0.19%"
    10 = "-1.42%"
    11 = "-0.29%"
    12 = "-2.85%"
    13 = "+2.32%"
    14 = "-0.60%"
    15 = "+1.54%"
    16 = "+1.98%"
    17 = "+2.09%"
    18 = "+1.00%"
    19 = "-0.82%"
    20 = "+3.28%"
    21 = "+0.36%"
    22 = "+0.87%"
    23 = "+0.44%"
    24 = "-0.51%"
    25 = "+0.22%"
    26 = "-0.04%"
    27 = "+0.18%"
    28 = "+7.19%"
    29 = "+1.99%"
    30 = "-2.25%"
    31 = "+0.09%"
    32 = "+1.34%"
    33 = "+5.79%"
    34 = "-1.81%"
    35 = "+4.53%"
    36 = "-2.23%"
    37 = "-0.82%"
    38 = "+0.74%"
    39 = "+0.76%"
    40 = "+17.67%"
    41 = "-0.18%"
    42 = "-28.13%"
    43 = "-1.47%"
    44 = "+1.64%"
    45 = "-4.56%"
    46 = "+1.76%"
    47 = "-0.95%"
    48 = "+1.47%"
    49 = "+0.42%"
    50 = "-3.53%"
    51 = "-0.06%"
}

# Rows whose new price text would otherwise be auto-parsed as a plain
# number by Excel (single decimal point) -- these must be forced to
# stay text so the stored value keeps its original textual form
# (e.g. "1.00" instead of becoming the number 1).
$forceTextRows = @(4, 5, 6, 10, 11, 12, 19, 20, 22, 23, 24, 26, 27, 28, 29, 30, 31, 33, 34, 36, 39, 40, 42, 45, 49, 51)

foreach ($row in $dValues.Keys) {
    $cell = $ws.Range("D$row")
    if ($forceTextRows -contains $row) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $dValues[$row]
}

foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = "  " + $eValues[$row] + "  "
}
